$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename task in B4 and update related status/date cells ---
$ws.Range("B4").Value = "Pencarian dan Analisis Dataset"

# --- Fill in "Tanggal Penyelesaian" completion dates for finished tasks ---
$ws.Range("E4").Value = [DateTime]::FromOADate(44491)
$ws.Range("E5").Value = [DateTime]::FromOADate(44493)
$ws.Range("E6").Value = [DateTime]::FromOADate(44495)

# --- Update status text for rows that are now finished / in progress ---
$ws.Range("G4").Value = "Done"
$ws.Range("G5").Value = "Done"
$ws.Range("G6").Value = "Done"
$ws.Range("G7").Value = "On Going"

# --- Color-code the Status column: Done = blue, On Going = green, Waiting = orange ---
$ws.Range("G4").Interior.Color = 15773696   # FF00B0F0 (blue)   -> Done
$ws.Range("G5").Interior.Color = 15773696   # FF00B0F0 (blue)   -> Done
$ws.Range("G6").Interior.Color = 15773696   # FF00B0F0 (blue)   -> Done
$ws.Range("G7").Interior.Color = 5287936    # FF92D050 (green)  -> On Going
$ws.Range("G8").Interior.Color = 3411545    # theme accent2 tint -> Waiting
$ws.Range("G9").Interior.Color = 3411545    # theme accent2 tint -> Waiting

# --- Bold the title and table header row ---
$ws.Range("D1").Font.Bold = $true
$ws.Range("A3:G3").Font.Bold = $true

# --- Widen JobDesc column to fit the longer text ---
$ws.Columns.Item(2).ColumnWidth = 27.6

# --- View settings: zoom + active selection ---
$excel.ActiveWindow.Zoom = 110
$ws.Range("K7").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1
